# Updated symbol list on Thu Jan 19 19:54:49 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# symbol table with the latest scraped quotes, and fixes the swapped
# Dexo / BKEXToken rows (42-43) so the coin name, link, price and volume
# line up correctly again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text such as
# "293.08" or "-0.23%". Force the range to Text format first so Excel
# doesn't silently reinterpret the assigned strings as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "293.08"
$ws.Range("E2").Value = "-0.23%"

$ws.Range("D3").Value = "30.93"
$ws.Range("E3").Value = "-0.42%"

$ws.Range("D4").Value = "4.902"
$ws.Range("E4").Value = "-0.10%"

$ws.Range("D5").Value = "0.07327"
$ws.Range("E5").Value = "0.24%"

$ws.Range("D6").Value = "2.308"
$ws.Range("E6").Value = "26.81%"

$ws.Range("D7").Value = "7.702"
$ws.Range("E7").Value = "0.12%"

$ws.Range("D8").Value = "3.729"
$ws.Range("E8").Value = "-0.91%"

$ws.Range("D9").Value = "0.9018"
$ws.Range("E9").Value = "-0.46%"

$ws.Range("D10").Value = "0.1681"
$ws.Range("E10").Value = "0.93%"

$ws.Range("D11").Value = "0.07940"
$ws.Range("E11").Value = "6.36%"

$ws.Range("D12").Value = "0.08160"
$ws.Range("E12").Value = "-0.07%"

$ws.Range("D13").Value = "0.03100"
$ws.Range("E13").Value = "3.61%"

$ws.Range("D14").Value = "0.1006"
$ws.Range("E14").Value = "0.57%"

$ws.Range("D15").Value = "0.001502"
$ws.Range("E15").Value = "-0.13%"

$ws.Range("D16").Value = "0.005685"
$ws.Range("E16").Value = "0.43%"

$ws.Range("D17").Value = "3.476"
$ws.Range("E17").Value = "0.40%"

$ws.Range("E18").Value = "-1.51%"

$ws.Range("D19").Value = "0.3330"
$ws.Range("E19").Value = "1.26%"

$ws.Range("D20").Value = "0.1300"
$ws.Range("E20").Value = "-0.40%"

$ws.Range("D21").Value = "4.016"
$ws.Range("E21").Value = "-7.11%"

$ws.Range("E22").Value = "4.89%"

$ws.Range("D23").Value = "0.04523"
$ws.Range("E23").Value = "0.90%"

$ws.Range("E24").Value = "-1.26%"

$ws.Range("D25").Value = "0.004641"
$ws.Range("E25").Value = "14.82%"

$ws.Range("E26").Value = "3.62%"

$ws.Range("D27").Value = "0.0003388"

$ws.Range("D39").Value = "0.01602"
$ws.Range("E39").Value = "-3.18%"

$ws.Range("D40").Value = "0.04439"
$ws.Range("E40").Value = "0.77%"

$ws.Range("D41").Value = "0.007345"
$ws.Range("E41").Value = "-1.52%"

# Rows 42/43 were swapped upstream: Dexo and BKEXToken traded places.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1325"
$ws.Range("E42").Value = "0.30%"

$ws.Range("B43").Value = "Dexo"
$ws.Range("C43").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D43").Value = "0.008621"
$ws.Range("E43").Value = "--%"

$ws.Range("D44").Value = "0.001998"
$ws.Range("E44").Value = "-2.13%"

$ws.Range("D45").Value = "0.009426"
$ws.Range("E45").Value = "-7.46%"

$ws.Range("D46").Value = "0.00005922"
$ws.Range("E46").Value = "-0.81%"

$ws.Range("E47").Value = "-0.35%"

$ws.Range("E48").Value = "3.83%"

$ws.Range("D49").Value = "0.002893"
$ws.Range("E49").Value = "18.89%"

$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").Value = "-0.35%"

$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").Value = "-0.35%"
